$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 4.17
$ws.Range("F2").Value = 164.56
$ws.Range("G2").Value = 10.51
$ws.Range("M2").Value = -0.34
$ws.Range("N2").Value = 32.02

# Row 3
$ws.Range("E3").Value = -2.15
$ws.Range("F3").Value = 167.73
$ws.Range("M3").Value = 1.38
$ws.Range("N3").Value = -0.46

# Row 4
$ws.Range("E4").Value = 205.05
$ws.Range("F4").Value = 1.52
$ws.Range("M4").Value = 0.47
$ws.Range("N4").Value = 0.35

# Row 5
$ws.Range("E5").Value = 9.57
$ws.Range("F5").Value = 3.58
$ws.Range("G5").Value = 1.51
$ws.Range("M5").Value = 1043.13
$ws.Range("N5").Value = 11.08
$ws.Range("O5").Value = 640.86

# Row 6
$ws.Range("E6").Value = 76.29000000000001
$ws.Range("F6").Value = 0.6899999999999999
$ws.Range("M6").Value = 2682.64
$ws.Range("N6").Value = 334.37
$ws.Range("O6").Value = 14227.79
$ws.Range("T6").Value = 1.1
